$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SIQ")

# v1.5 review pass: reword the Admin/Login/Registration questions per the
# review comments on v1.2 (login, registration and ID constraints).
$ws.Range("C18").Value = "Should the admin be able to remove a user from the system?"
$ws.Range("C20").Value = "Should access to the system be restricted to users who have registered?"
$ws.Range("C21").Value = "Should the system display a simple error message if login information is incorrect or incomplete?"
$ws.Range("C22").Value = "Should user passwords be stored securely using basic hashing methods? "
$ws.Range("C23").Value = "Should the registration form require all fields (email, username, and password) to be filled out?"
$ws.Range("C24").Value = "Should the username have specific constraints (e.g., 4-20 characters,only contain letters and numbers  no special symbols like(!,@,#,`$<%,^,etc))"
$ws.Range("C25").Value = "Should the system prevent multiple registrations using the same email or username?"
$ws.Range("C26").Value = " Should the system automatically assign user IDs in a simple, consistent format (e.g., U001, U002)?"

# Switch the active/selected tab back to the SIQ sheet (was "Version history").
$ws.Activate()
